# Updates the cryptos list with the latest scraped prices / volume figures.
# (Mirrors the automated "Updated cryptos list ... with GitHub Actions" commit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '37.059.04'
    'E2' = '  -1.02%  '
    'D3' = '1.992.95'
    'E3' = '  -2.31%  '
    'D4' = '1.00'
    'E4' = '  -0.19%  '
    'D5' = '270.01'
    'E5' = '  +8.84%  '
    'D6' = '0.610'
    'E6' = '  -2.16%  '
    'E7' = '  +0.15%  '
    'D8' = '55.65'
    'E8' = '  -6.00%  '
    'D9' = '0.377'
    'E9' = '  -4.26%  '
    'D10' = '0.0764'
    'E10' = '  -5.33%  '
    'E11' = '  -3.36%  '
    'D12' = '14.35'
    'E12' = '  -5.35%  '
    'D13' = '2.256.36'
    'E13' = '  -3.58%  '
    'D14' = '21.61'
    'E14' = '  -2.07%  '
    'D15' = '0.773'
    'E15' = '  -9.09%  '
    'D16' = '5.18'
    'E16' = '  -4.95%  '
    'D17' = '2.004.21'
    'E17' = '  -1.77%  '
    'D18' = '36.926.84'
    'E18' = '  -1.19%  '
    'D19' = '69.51'
    'E19' = '  -1.31%  '
    'D20' = '0.0₃0828'
    'E20' = '  -3.99%  '
    'D21' = '234.37'
    'E21' = '  +1.94%  '
    'D22' = '5.06'
    'E22' = '  -4.20%  '
    'B23' = 'PancakeSwap'
    'C23' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D23' = '2.66'
    'E23' = '  +4.08%  '
    'B24' = 'Dai'
    'C24' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D24' = '1.00'
    'E24' = '  +0.13%  '
    'D25' = '2.38'
    'E25' = '  +0.87%  '
    'D26' = '164.22'
    'E26' = '  -0.21%  '
    'D27' = '8.90'
    'E27' = '  -4.73%  '
    'D28' = '19.33'
    'E28' = '  -3.03%  '
    'B29' = 'Kaspa'
    'C29' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D29' = '0.125'
    'E29' = '  -9.72%  '
    'B30' = 'ImmutableX'
    'C30' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D30' = '1.34'
    'E30' = '  -2.28%  '
    'E31' = '  -2.63%  '
    'D32' = '4.53'
    'E32' = '  -6.02%  '
    'D33' = '0.0625'
    'E33' = '  -8.39%  '
    'D34' = '4.35'
    'E34' = '  -3.59%  '
    'D35' = '2.39'
    'E35' = '  -5.36%  '
    'D36' = '3.45'
    'E36' = '  -4.21%  '
    'E37' = '  -0.27%  '
    'D38' = '0.999'
    'E38' = '  -0.31%  '
    'D39' = '5.36'
    'E39' = '  -2.08%  '
    'D40' = '3.00'
    'E40' = '  -0.21%  '
    'B41' = 'TrustWalletToken'
    'C41' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D41' = '1.18'
    'E41' = '  -0.17%  '
    'B42' = 'Maker'
    'C42' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D42' = '1.444.90'
    'E42' = '  +2.88%  '
    'D43' = '0.0913'
    'E43' = '  -7.06%  '
    'D44' = '0.0207'
    'E44' = '  -4.57%  '
    'D45' = '88.65'
    'E45' = '  -3.48%  '
    'D46' = '15.54'
    'E46' = '  -7.18%  '
    'D47' = '1.02'
    'E47' = '  -4.19%  '
    'D48' = '2.91'
    'E48' = '  +0.65%  '
    'D49' = '6.81'
    'E49' = '  -9.29%  '
    'B50' = 'NEARProtocol'
    'C50' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D50' = '1.93'
    'E50' = '  -8.33%  '
    'B51' = 'RocketPoolETH'
    'C51' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D51' = '2.148.49'
    'E51' = '  -3.63%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the cell to be treated as literal text (these columns store
    # formatted/rounded numbers and percentages as strings, not numbers),
    # then restore the original (default/"Normal") style so no formatting
    # is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
